$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-45: column A = Actual Consumption (MW), column B = Timestamp (serial date)
$data = @(
    ,@(2, 5618, 45741)
    ,@(3, 5550, 45741.01041666666)
    ,@(4, 5514, 45741.02083333334)
    ,@(5, 5488, 45741.03125)
    ,@(6, 5447, 45741.04166666666)
    ,@(7, 5414, 45741.05208333334)
    ,@(8, 5411, 45741.0625)
    ,@(9, 5415, 45741.07291666666)
    ,@(10, 5330, 45741.08333333334)
    ,@(11, 5310, 45741.09375)
    ,@(12, 5307, 45741.10416666666)
    ,@(13, 5318, 45741.11458333334)
    ,@(14, 5339, 45741.125)
    ,@(15, 5370, 45741.13541666666)
    ,@(16, 5422, 45741.14583333334)
    ,@(17, 5464, 45741.15625)
    ,@(18, 5528, 45741.16666666666)
    ,@(19, 5592, 45741.17708333334)
    ,@(20, 5669, 45741.1875)
    ,@(21, 5744, 45741.19791666666)
    ,@(22, 5856, 45741.20833333334)
    ,@(23, 5917, 45741.21875)
    ,@(24, 6039, 45741.22916666666)
    ,@(25, 6176, 45741.23958333334)
    ,@(26, 6367, 45741.25)
    ,@(27, 6444, 45741.26041666666)
    ,@(28, 6526, 45741.27083333334)
    ,@(29, 6537, 45741.28125)
    ,@(30, 6608, 45741.29166666666)
    ,@(31, 6527, 45741.30208333334)
    ,@(32, 6528, 45741.3125)
    ,@(33, 6475, 45741.32291666666)
    ,@(34, 6330, 45741.33333333334)
    ,@(35, 6286, 45741.34375)
    ,@(36, 6184, 45741.35416666666)
    ,@(37, 6101, 45741.36458333334)
    ,@(38, 5946, 45741.375)
    ,@(39, 5894, 45741.38541666666)
    ,@(40, 5772, 45741.39583333334)
    ,@(41, 5729, 45741.40625)
    ,@(42, 5592, 45741.41666666666)
    ,@(43, 5555, 45741.42708333334)
    ,@(44, 5515, 45741.4375)
    ,@(45, 5564, 45741.44791666666)
)

foreach ($row in $data) {
    $r = $row[0]
    $aVal = $row[1]
    $bVal = $row[2]
    $ws.Cells.Item($r, 1).Value = $aVal
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $bVal
    if ($r -gt 33) {
        $cellB.NumberFormat = $ws.Cells.Item(33, 2).NumberFormat
    }
}